# Append the 8 new NBA game rows (rows 649-656) to Sheet1, matching the
# target diff: new shared string "Yes", new data rows with the same
# B/D (points) number-format style as the rows immediately above them,
# and a sheet view reset back to the top of the sheet (no stale
# scroll/selection left pointing at the old last row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Charlotte Hornets",       106, "Detroit Pistons",        113, "No",  17832, "Little Caesars Arena",     "Detroit Pistons",        "Charlotte Hornets"),
    @("Minnesota Timberwolves",  118, "Washington Wizards",     107, "No",  17832, "Capital One Arena",        "Minnesota Timberwolves", "Washington Wizards"),
    @("Memphis Grizzlies",       105, "Miami Heat",               96, "No",  17832, "Kaseya Center",            "Memphis Grizzlies",      "Miami Heat"),
    @("Portland Trail Blazers",  137, "Houston Rockets",        131, "Yes", 17832, "Toyota Center",             "Portland Trail Blazers", "Houston Rockets"),
    @("Cleveland Cavaliers",     116, "Milwaukee Bucks",        126, "No",  17832, "Fiserv Forum",              "Milwaukee Bucks",        "Cleveland Cavaliers"),
    @("Phoenix Suns",            132, "Dallas Mavericks",       109, "No",  17832, "American Airlines Center", "Phoenix Suns",           "Dallas Mavericks"),
    @("Oklahoma City Thunder",   140, "San Antonio Spurs",      114, "No",  17832, "Frost Bank Center",         "Oklahoma City Thunder",  "San Antonio Spurs"),
    @("Atlanta Hawks",           112, "Golden State Warriors",  134, "No",  17832, "Chase Center",              "Golden State Warriors",  "Atlanta Hawks")
)

$startRow = 649
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $arr = New-Object 'object[,]' 1,9
    for ($c = 0; $c -lt 9; $c++) {
        $arr[0,$c] = $data[$c]
    }

    $ws.Range($ws.Cells.Item($r,1), $ws.Cells.Item($r,9)).Value = $arr

    # Away Pts / Home Pts keep the "#,##0" number style already used by the
    # preceding rows (cellXfs index 3 in styles.xml).
    $ws.Range($ws.Cells.Item($r,2), $ws.Cells.Item($r,2)).NumberFormat = "#,##0"
    $ws.Range($ws.Cells.Item($r,4), $ws.Cells.Item($r,4)).NumberFormat = "#,##0"
}

# Reset the view: scroll/select back to the top-left instead of leaving it
# parked on the old final row.
$ws.Range("A1").Select()
